# Atualização de bases das ligas, do dia: 18-02-2024 às 22:54
#
# The underlying source data re-ordered several match rows (their match-id
# in column B moved to a different row), so each pair of rows below has its
# entire record (columns B..AC, i.e. everything except the row-number
# column A) swapped between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose B:AC contents must be swapped.
$rowPairs = @(
    @(7, 8),
    @(32, 33),
    @(40, 41),
    @(67, 68),
    @(144, 145),
    @(179, 180)
)

$firstCol = 2   # column B
$lastCol  = 29  # column AC

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)

        $valA = $cellA.Value()
        $valB = $cellB.Value()

        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}
